# Applies the "plotted toy-spam with min 5" edit: refreshed anchor-word
# statistics across the sheet, re-ordered the shared-string table to match
# the regenerated word list, and appended two new anchor rows (30, 31).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: wipe the whole existing table so every shared string currently
# in xl/sharedStrings.xml becomes unreferenced and gets dropped. Re-populating
# the cells afterwards (in the order below) rebuilds the shared-string table
# from index 0, in exactly the order the target workbook uses. ---
$ws.Range("A1:Q29").ClearContents()

# --- Step 2: re-enter every string cell in the precise order the target
# file lists its strings in (indices 0..50 of xl/sharedStrings.xml) ---
$ws.Range("A2").Value = "name"
$ws.Range("B2").Value = "anchor score"
$ws.Range("C2").Value = "type occurences"
$ws.Range("D2").Value = "total occurences"
$ws.Range("E2").Value = "+%"
$ws.Range("F2").Value = "-%"
$ws.Range("G2").Value = "both"
$ws.Range("H2").Value = "normal"
$ws.Range("J2").Value = "name"
$ws.Range("K2").Value = "anchor score"
$ws.Range("L2").Value = "type occurences"
$ws.Range("M2").Value = "total occurences"
$ws.Range("N2").Value = "+%"
$ws.Range("O2").Value = "-%"
$ws.Range("P2").Value = "both"
$ws.Range("Q2").Value = "normal"
$ws.Range("A3").Value = "poorly"
$ws.Range("A4").Value = "disappointing"
$ws.Range("A5").Value = "poor"
$ws.Range("A6").Value = "disappointed"
$ws.Range("A7").Value = "however"
$ws.Range("A8").Value = "broke"
$ws.Range("A9").Value = "waste"
$ws.Range("A10").Value = "junk"
$ws.Range("A11").Value = "smaller"
$ws.Range("A12").Value = "small"
$ws.Range("A13").Value = "plastic"
$ws.Range("A14").Value = "broken"
$ws.Range("A15").Value = "apart"
$ws.Range("A16").Value = "difficult"
$ws.Range("A17").Value = "ok"
$ws.Range("A18").Value = "cheap"
$ws.Range("A19").Value = "thought"
$ws.Range("A20").Value = "item"
$ws.Range("A21").Value = "size"
$ws.Range("A22").Value = "could"
$ws.Range("A23").Value = "used"
$ws.Range("A24").Value = "would"
$ws.Range("A25").Value = "money"
$ws.Range("A26").Value = "work"
$ws.Range("A27").Value = "better"
$ws.Range("A28").Value = "product"
$ws.Range("A29").Value = "price"
$ws.Range("A30").Value = "use"
$ws.Range("A31").Value = "like"
$ws.Range("A1").Value = "negative"
$ws.Range("J3").Value = "awesome"
$ws.Range("J4").Value = "wonderful"
$ws.Range("J5").Value = "favorite"
$ws.Range("J6").Value = "excellent"
$ws.Range("J7").Value = "thank"
$ws.Range("J8").Value = "great"
$ws.Range("J9").Value = "love"
$ws.Range("J10").Value = "loves"
$ws.Range("J11").Value = "perfect"
$ws.Range("J12").Value = "loved"
$ws.Range("J13").Value = "fun"
$ws.Range("J14").Value = "game"
$ws.Range("J1").Value = "positive"

# --- Step 3: restore the bold/border/center-top header style on row 2 and
# the A-column anchor labels (style index 1 in xl/styles.xml) ---
$ws.Range("A2:H2,J2:Q2").Font.Bold = $true
$ws.Range("A2:H2,J2:Q2").HorizontalAlignment = -4108
$ws.Range("A2:H2,J2:Q2").VerticalAlignment = -4160
$ws.Range("A2:H2,J2:Q2").Borders.LineStyle = 1
$ws.Range("A3:A29").Font.Bold = $true
$ws.Range("A3:A29").HorizontalAlignment = -4108
$ws.Range("A3:A29").VerticalAlignment = -4160
$ws.Range("A3:A29").Borders.LineStyle = 1
$ws.Range("J3:J14").Font.Bold = $true
$ws.Range("J3:J14").HorizontalAlignment = -4108
$ws.Range("J3:J14").VerticalAlignment = -4160
$ws.Range("J3:J14").Borders.LineStyle = 1

# --- Step 4: update the numeric / boolean statistics for rows 3-29 ---
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 46
$ws.Range("D3").Value = 46
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = $False
$ws.Range("H3").Value = 0
$ws.Range("K3").Value = 0.8461538461538461
$ws.Range("L3").Value = 55
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = $False
$ws.Range("Q3").Value = 10
$ws.Range("B4").Value = 0.8409090909090909
$ws.Range("C4").Value = 37
$ws.Range("D4").Value = 37
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = $False
$ws.Range("H4").Value = 7
$ws.Range("K4").Value = 0.8392857142857143
$ws.Range("L4").Value = 47
$ws.Range("M4").Value = 47
$ws.Range("N4").Value = 1
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = $False
$ws.Range("Q4").Value = 9
$ws.Range("B5").Value = 0.7464788732394366
$ws.Range("C5").Value = 53
$ws.Range("D5").Value = 53
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = $False
$ws.Range("H5").Value = 18
$ws.Range("K5").Value = 0.6989247311827957
$ws.Range("L5").Value = 65
$ws.Range("M5").Value = 65
$ws.Range("N5").Value = 1
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $False
$ws.Range("Q5").Value = 28
$ws.Range("B6").Value = 0.7365591397849462
$ws.Range("C6").Value = 137
$ws.Range("D6").Value = 137
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = $False
$ws.Range("H6").Value = 49
$ws.Range("K6").Value = 0.46875
$ws.Range("L6").Value = 30
$ws.Range("M6").Value = 30
$ws.Range("N6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = $False
$ws.Range("Q6").Value = 34
$ws.Range("B7").Value = 0.71875
$ws.Range("C7").Value = 46
$ws.Range("D7").Value = 46
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = $False
$ws.Range("H7").Value = 18
$ws.Range("K7").Value = 0.4347826086956522
$ws.Range("L7").Value = 30
$ws.Range("M7").Value = 30
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = $False
$ws.Range("Q7").Value = 39
$ws.Range("B8").Value = 0.7038834951456311
$ws.Range("C8").Value = 145
$ws.Range("D8").Value = 145
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = $False
$ws.Range("H8").Value = 61
$ws.Range("K8").Value = 0.3459016393442623
$ws.Range("L8").Value = 422
$ws.Range("M8").Value = 422
$ws.Range("N8").Value = 1
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = $False
$ws.Range("Q8").Value = 798
$ws.Range("B9").Value = 0.6283783783783784
$ws.Range("C9").Value = 93
$ws.Range("D9").Value = 93
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = $False
$ws.Range("H9").Value = 55
$ws.Range("K9").Value = 0.2998565279770445
$ws.Range("L9").Value = 209
$ws.Range("M9").Value = 209
$ws.Range("N9").Value = 1
$ws.Range("O9").Value = 0
$ws.Range("P9").Value = $False
$ws.Range("Q9").Value = 488
$ws.Range("B10").Value = 0.6
$ws.Range("C10").Value = 33
$ws.Range("D10").Value = 33
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = $False
$ws.Range("H10").Value = 22
$ws.Range("K10").Value = 0.2489626556016598
$ws.Range("L10").Value = 120
$ws.Range("M10").Value = 120
$ws.Range("N10").Value = 1
$ws.Range("O10").Value = 0
$ws.Range("P10").Value = $False
$ws.Range("Q10").Value = 362
$ws.Range("B11").Value = 0.5882352941176471
$ws.Range("C11").Value = 70
$ws.Range("D11").Value = 70
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = $False
$ws.Range("H11").Value = 49
$ws.Range("K11").Value = 0.2168674698795181
$ws.Range("L11").Value = 36
$ws.Range("M11").Value = 36
$ws.Range("N11").Value = 1
$ws.Range("O11").Value = 0
$ws.Range("P11").Value = $False
$ws.Range("Q11").Value = 130
$ws.Range("B12").Value = 0.4753623188405797
$ws.Range("C12").Value = 164
$ws.Range("D12").Value = 164
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = $False
$ws.Range("H12").Value = 181
$ws.Range("K12").Value = 0.1896024464831804
$ws.Range("L12").Value = 62
$ws.Range("M12").Value = 62
$ws.Range("N12").Value = 1
$ws.Range("O12").Value = 0
$ws.Range("P12").Value = $False
$ws.Range("Q12").Value = 265
$ws.Range("B13").Value = 0.4173228346456693
$ws.Range("C13").Value = 53
$ws.Range("D13").Value = 53
$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = $False
$ws.Range("H13").Value = 74
$ws.Range("K13").Value = 0.08413672217353199
$ws.Range("L13").Value = 96
$ws.Range("M13").Value = 96
$ws.Range("N13").Value = 1
$ws.Range("O13").Value = 0
$ws.Range("P13").Value = $False
$ws.Range("Q13").Value = 1045
$ws.Range("B14").Value = 0.4096385542168675
$ws.Range("C14").Value = 34
$ws.Range("D14").Value = 34
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = $False
$ws.Range("H14").Value = 49
$ws.Range("K14").Value = 0.03311688311688311
$ws.Range("L14").Value = 51
$ws.Range("M14").Value = 52
$ws.Range("N14").Value = 0.98
$ws.Range("O14").Value = 0.02000000000000002
$ws.Range("P14").Value = $True
$ws.Range("Q14").Value = 1489
$ws.Range("B15").Value = 0.3894736842105263
$ws.Range("C15").Value = 37
$ws.Range("D15").Value = 37
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = $False
$ws.Range("H15").Value = 58
$ws.Range("B16").Value = 0.3370786516853932
$ws.Range("C16").Value = 30
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = $False
$ws.Range("H16").Value = 59
$ws.Range("B17").Value = 0.328125
$ws.Range("C17").Value = 42
$ws.Range("D17").Value = 42
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = $False
$ws.Range("H17").Value = 86
$ws.Range("B18").Value = 0.2796208530805687
$ws.Range("C18").Value = 59
$ws.Range("D18").Value = 59
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 1
$ws.Range("G18").Value = $False
$ws.Range("H18").Value = 152
$ws.Range("B19").Value = 0.2722772277227723
$ws.Range("C19").Value = 55
$ws.Range("D19").Value = 55
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 1
$ws.Range("G19").Value = $False
$ws.Range("H19").Value = 147
$ws.Range("B20").Value = 0.213768115942029
$ws.Range("C20").Value = 59
$ws.Range("D20").Value = 59
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 1
$ws.Range("G20").Value = $False
$ws.Range("H20").Value = 217
$ws.Range("B21").Value = 0.1958762886597938
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 38
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 1
$ws.Range("G21").Value = $False
$ws.Range("H21").Value = 156
$ws.Range("B22").Value = 0.1847133757961783
$ws.Range("C22").Value = 29
$ws.Range("D22").Value = 29
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = $False
$ws.Range("H22").Value = 128
$ws.Range("B23").Value = 0.1828571428571429
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 32
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = $False
$ws.Range("H23").Value = 143
$ws.Range("B24").Value = 0.1810089020771513
$ws.Range("C24").Value = 122
$ws.Range("D24").Value = 122
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = $False
$ws.Range("H24").Value = 552
$ws.Range("B25").Value = 0.180379746835443
$ws.Range("C25").Value = 57
$ws.Range("D25").Value = 57
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = $False
$ws.Range("H25").Value = 259
$ws.Range("B26").Value = 0.1582278481012658
$ws.Range("C26").Value = 50
$ws.Range("D26").Value = 50
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = 1
$ws.Range("G26").Value = $False
$ws.Range("H26").Value = 266
$ws.Range("B27").Value = 0.1448598130841121
$ws.Range("C27").Value = 31
$ws.Range("D27").Value = 31
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = $False
$ws.Range("H27").Value = 183
$ws.Range("B28").Value = 0.1387665198237885
$ws.Range("C28").Value = 63
$ws.Range("D28").Value = 63
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 1
$ws.Range("G28").Value = $False
$ws.Range("H28").Value = 391
$ws.Range("B29").Value = 0.1005747126436782
$ws.Range("C29").Value = 35
$ws.Range("D29").Value = 35
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = 1
$ws.Range("G29").Value = $False
$ws.Range("H29").Value = 313

# --- Step 5: append the two new anchor rows, copying A29's header style
# (bold + border + center/top alignment = style index 1) onto the new labels ---
$ws.Range("B30").Value = 0.07945205479452055
$ws.Range("C30").Value = 29
$ws.Range("D30").Value = 29
$ws.Range("E30").Value = 0
$ws.Range("F30").Value = 1
$ws.Range("G30").Value = $False
$ws.Range("H30").Value = 336
$ws.Range("A29").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("B31").Value = 0.06085526315789474
$ws.Range("C31").Value = 37
$ws.Range("D31").Value = 37
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 1
$ws.Range("G31").Value = $False
$ws.Range("H31").Value = 571
$ws.Range("A29").Copy()
$ws.Range("A31").PasteSpecial(-4122)
